$d = $word.ActiveDocument

# Locate the paragraph that ends with "...documentações de software, etc"
# (the last sentence of section 4.2. Infraestrutura).
$searchText = "documentações de software, etc"
$r = $d.Content
$found = $r.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Figure out which paragraph index the match falls into, then grab the
    # paragraph right after it -- this is the existing blank, centered
    # paragraph that follows the sentence.
    $upTo = $d.Range(0, $r.End)
    $paraIndex = $upTo.Paragraphs.Count
    $followingPara = $d.Paragraphs.Item($paraIndex + 1)

    # Insert two new blank paragraphs before it, inheriting its formatting
    # (centered, spacing 360/auto, sz 24) -- matching the two new empty
    # centered paragraphs added by the edit.
    $followingPara.Range.InsertParagraphBefore()
    $followingPara.Range.InsertParagraphBefore()
}
